$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 4907087.5
$ws.Range("I76").Value = 12824806
$ws.Range("J76").Value = 5642.857
$ws.Range("K76").Value = 12824806
$ws.Range("L76").Value = 5642.857
$ws.Range("M76").Value = -12824491
$ws.Range("N76").Value = -6272.857
$ws.Range("H79").Value = 4907087.5
$ws.Range("I79").Value = 12824806
$ws.Range("J79").Value = 5642.857
$ws.Range("K79").Value = 12824806
$ws.Range("L79").Value = 5642.857
$ws.Range("M79").Value = -12823714
$ws.Range("N79").Value = -7826.857
$ws.Range("H80").Value = 6321.2354
$ws.Range("I80").Value = 380
$ws.Range("J80").Value = 20580.2
$ws.Range("K80").Value = 1140
$ws.Range("L80").Value = 61740.60000000001
$ws.Range("M80").Value = -142
$ws.Range("N80").Value = -63736.60000000001
$ws.Range("H83").Value = 6321.2354
$ws.Range("I83").Value = 380
$ws.Range("J83").Value = 20580.2
$ws.Range("K83").Value = 3420
$ws.Range("L83").Value = 185221.8
$ws.Range("M83").Value = 1572
$ws.Range("N83").Value = -195205.8
$ws.Range("H113").Value = 12502161
$ws.Range("I113").Value = 2229.2856
$ws.Range("J113").Value = 41668668
$ws.Range("K113").Value = 2229.2856
$ws.Range("L113").Value = 41668668
$ws.Range("M113").Value = 1024.7144
$ws.Range("N113").Value = -41675176
$ws.Range("H128").Value = 79800
$ws.Range("J128").Value = 79800
$ws.Range("L128").Value = 79800
$ws.Range("N128").Value = -89760
$ws.Range("H138").Value = 4336.145
$ws.Range("I138").Value = 2305.2
$ws.Range("J138").Value = 4900.2964
$ws.Range("K138").Value = 6915.599999999999
$ws.Range("L138").Value = 14700.8892
$ws.Range("M138").Value = -1775.599999999999
$ws.Range("N138").Value = -24980.8892

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 20084.191
$ws.Range("I32").Value = 17852
$ws.Range("K32").Value = 17852
$ws.Range("M32").Value = -17565
$ws.Range("H61").Value = 351611.25
$ws.Range("I61").Value = 8801.823
$ws.Range("J61").Value = 837257.9399999999
$ws.Range("K61").Value = 8801.823
$ws.Range("L61").Value = 837257.9399999999
$ws.Range("M61").Value = -8589.823
$ws.Range("N61").Value = -837681.9399999999
$ws.Range("H63").Value = 50011170
$ws.Range("I63").Value = 100002340
$ws.Range("J63").Value = 20000
$ws.Range("K63").Value = 100002340
$ws.Range("L63").Value = 20000
$ws.Range("M63").Value = -100001654
$ws.Range("N63").Value = -21372
$ws.Range("H66").Value = 50011170
$ws.Range("I66").Value = 100002340
$ws.Range("J66").Value = 20000
$ws.Range("K66").Value = 500011700
$ws.Range("L66").Value = 100000
$ws.Range("M66").Value = -500008268
$ws.Range("N66").Value = -106864
$ws.Range("H88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("L88").ClearContents()
$ws.Range("N88").Value = 0
$ws.Range("H91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("L91").ClearContents()
$ws.Range("N91").Value = 0
$ws.Range("H102").Value = 1852998.4
$ws.Range("I102").Value = 1852998.4
$ws.Range("K102").Value = 1852998.4
$ws.Range("M102").Value = -1851376.4
$ws.Range("H136").Value = 351611.25
$ws.Range("I136").Value = 8801.823
$ws.Range("J136").Value = 837257.9399999999
$ws.Range("K136").Value = 26405.469
$ws.Range("L136").Value = 2511773.82
$ws.Range("M136").Value = -23855.469
$ws.Range("N136").Value = -2516873.82

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 4551.6665
$ws.Range("I105").Value = 2166.3333
$ws.Range("J105").Value = 6937
$ws.Range("K105").Value = 2166.3333
$ws.Range("L105").Value = 6937
$ws.Range("M105").Value = -419.3332999999998
$ws.Range("N105").Value = -10431
$ws.Range("H134").Value = 25642.543
$ws.Range("I134").Value = 4159.5854
$ws.Range("J134").Value = 201802.8
$ws.Range("K134").Value = 12478.7562
$ws.Range("L134").Value = 605408.3999999999
$ws.Range("M134").Value = -9943.7562
$ws.Range("N134").Value = -610478.3999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 16133109
$ws.Range("I31").Value = 1403.0476
$ws.Range("J31").Value = 50009692
$ws.Range("K31").Value = 1403.0476
$ws.Range("L31").Value = 50009692
$ws.Range("M31").Value = -1108.0476
$ws.Range("N31").Value = -50010282
$ws.Range("H34").Value = 16133109
$ws.Range("I34").Value = 1403.0476
$ws.Range("J34").Value = 50009692
$ws.Range("K34").Value = 1403.0476
$ws.Range("L34").Value = 50009692
$ws.Range("M34").Value = -1201.0476
$ws.Range("N34").Value = -50010096
$ws.Range("H62").Value = 4854.6
$ws.Range("I62").Value = 5170.6924
$ws.Range("K62").Value = 5170.6924
$ws.Range("M62").Value = -4546.6924
$ws.Range("H65").Value = 4854.6
$ws.Range("I65").Value = 5170.6924
$ws.Range("K65").Value = 25853.462
$ws.Range("M65").Value = -22733.462
$ws.Range("H99").Value = 2600
$ws.Range("I99").Value = 1256
$ws.Range("J99").Value = 3944
$ws.Range("K99").Value = 1256
$ws.Range("L99").Value = 3944
$ws.Range("M99").Value = 242
$ws.Range("N99").Value = -6940
$ws.Range("H107").Value = 1463.381
$ws.Range("I107").Value = 927
$ws.Range("K107").Value = 927
$ws.Range("M107").Value = 993
$ws.Range("H126").Value = 2600
$ws.Range("I126").Value = 1256
$ws.Range("J126").Value = 3944
$ws.Range("K126").Value = 3768
$ws.Range("L126").Value = 11832
$ws.Range("M126").Value = -1298
$ws.Range("N126").Value = -16772
$ws.Range("H132").Value = 5129711
$ws.Range("I132").Value = 6061629
$ws.Range("K132").Value = 18184887
$ws.Range("M132").Value = -18182357
$ws.Range("H134").Value = 13737018
$ws.Range("I134").Value = 14496692
$ws.Range("J134").Value = 5000757
$ws.Range("K134").Value = 43490076
$ws.Range("L134").Value = 15002271
$ws.Range("M134").Value = -43487541
$ws.Range("N134").Value = -15007341
$ws.Range("H140").Value = 39085.184
$ws.Range("J140").Value = 39085.184
$ws.Range("L140").Value = 39085.184
$ws.Range("N140").Value = -49445.184

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5093.364
$ws.Range("I70").Value = 5007.143
$ws.Range("J70").Value = 5156.8945
$ws.Range("K70").Value = 5007.143
$ws.Range("L70").Value = 5156.8945
$ws.Range("M70").Value = -4737.143
$ws.Range("N70").Value = -5696.8945
$ws.Range("H73").Value = 5093.364
$ws.Range("I73").Value = 5007.143
$ws.Range("J73").Value = 5156.8945
$ws.Range("K73").Value = 5007.143
$ws.Range("L73").Value = 5156.8945
$ws.Range("M73").Value = -4071.143
$ws.Range("N73").Value = -7028.8945
$ws.Range("H80").Value = 10600.917
$ws.Range("I80").Value = 51202.5
$ws.Range("J80").Value = 2480.6
$ws.Range("K80").Value = 51202.5
$ws.Range("L80").Value = 2480.6
$ws.Range("M80").Value = -50204.5
$ws.Range("N80").Value = -4476.6
$ws.Range("H83").Value = 10600.917
$ws.Range("I83").Value = 51202.5
$ws.Range("J83").Value = 2480.6
$ws.Range("K83").Value = 256012.5
$ws.Range("L83").Value = 12403
$ws.Range("M83").Value = -251020.5
$ws.Range("N83").Value = -22387
$ws.Range("H102").Value = 3559
$ws.Range("I102").Value = 3000
$ws.Range("K102").Value = 3000
$ws.Range("M102").Value = -1378
$ws.Range("H132").Value = 7578105
$ws.Range("I132").Value = 9261128
$ws.Range("J132").Value = 4499.75
$ws.Range("K132").Value = 27783384
$ws.Range("L132").Value = 13499.25
$ws.Range("M132").Value = -27780854
$ws.Range("N132").Value = -18559.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 2750829.2
$ws.Range("I122").Value = 3250321.2
$ws.Range("K122").Value = 9750963.600000001
$ws.Range("M122").Value = -9748513.600000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1657.7142
$ws.Range("I122").Value = 1352
$ws.Range("J122").Value = 1780
$ws.Range("K122").Value = 4056
$ws.Range("L122").Value = 5340
$ws.Range("M122").Value = -1606
$ws.Range("N122").Value = -10240
